$wb = $excel.ActiveWorkbook

# --- Pilot 2 sheet: clear frozen-pane topLeftCell override & reset selection ---
$wsPilot2 = $wb.Worksheets.Item("Pilot 2")
$wsPilot2.Range("B2").Select()

# --- Quality Ratings sheet: add visual_cue column + recode values ---
$ws = $wb.Worksheets.Item("Quality Ratings")

$ws.Range("B2").Value2 = "a flower"
$ws.Range("C2").Value2 = "a flower"

$ws.Range("D1").Value2 = "visual_cue"
$ws.Range("D2").Value2 = "no"

$ws.Range("B3").Value2 = "a flower"
$ws.Range("C3").Value2 = "a flower"
$ws.Range("D3").Value2 = "no"

$ws.Range("B4").Value2 = "a flower"
$ws.Range("C4").Value2 = "a flower"
$ws.Range("D4").Value2 = "no"

$ws.Range("B5").Value2 = "a flower"
$ws.Range("C5").Value2 = "a flower"
$ws.Range("D5").Value2 = "no"

$ws.Range("B6").Value2 = "a flower"
$ws.Range("C6").Value2 = "a flower"
$ws.Range("D6").Value2 = "no"

$ws.Range("B7").Value2 = "a flower"
$ws.Range("C7").Value2 = "a flower"
$ws.Range("D7").Value2 = "no"

$ws.Range("B8").Value2 = "a flower"
$ws.Range("C8").Value2 = "a flower"
$ws.Range("D8").Value2 = "no"

$ws.Range("B9").Value2 = "a flower"
$ws.Range("C9").Value2 = "a flower"
$ws.Range("D9").Value2 = "no"

$ws.Range("B10").Value2 = "a flower"
$ws.Range("C10").Value2 = "a flower"
$ws.Range("D10").Value2 = "no"

$ws.Range("B11").Value2 = "a flower"
$ws.Range("C11").Value2 = "a flower"
$ws.Range("D11").Value2 = "no"

$ws.Range("B12").Value2 = "a flower"
$ws.Range("C12").Value2 = "a flower"
$ws.Range("D12").Value2 = "no"

$ws.Range("B13").Value2 = "a flower"
$ws.Range("C13").Value2 = "a flower"
$ws.Range("D13").Value2 = "no"

$ws.Range("B14").Value2 = "a flower garden"
$ws.Range("C14").Value2 = "a flower garden"
$ws.Range("D14").Value2 = "no"

$ws.Range("B15").Value2 = "a flower garden"
$ws.Range("C15").Value2 = "a flower garden"
$ws.Range("D15").Value2 = "no"

$ws.Range("B16").Value2 = "a flower garden"
$ws.Range("C16").Value2 = "a flower garden"
$ws.Range("D16").Value2 = "no"

$ws.Range("B17").Value2 = "a flower garden"
$ws.Range("C17").Value2 = "a flower garden"
$ws.Range("D17").Value2 = "no"

$ws.Range("B18").Value2 = "a flower garden"
$ws.Range("C18").Value2 = "a flower garden"
$ws.Range("D18").Value2 = "yes"

$ws.Range("B19").Value2 = "a flower garden"
$ws.Range("C19").Value2 = "a flower garden"
$ws.Range("D19").Value2 = "yes"

$ws.Range("B20").Value2 = "some flowers"
$ws.Range("C20").Value2 = "some flowers"
$ws.Range("D20").Value2 = "yes"

$ws.Range("B21").Value2 = "some flowers"
$ws.Range("C21").Value2 = "some flowers"
$ws.Range("D21").Value2 = "yes"

$ws.Range("B22").Value2 = "some flowers"
$ws.Range("C22").Value2 = "some flowers"
$ws.Range("D22").Value2 = "yes"

$ws.Range("B23").Value2 = "some flowers"
$ws.Range("C23").Value2 = "some flowers"
$ws.Range("D23").Value2 = "yes"

$ws.Range("B24").Value2 = "some flowers"
$ws.Range("C24").Value2 = "some flowers"
$ws.Range("D24").Value2 = "yes"

$ws.Range("B25").Value2 = "a slice of watermelon"
$ws.Range("C25").Value2 = "a slice of watermelon"
$ws.Range("D25").Value2 = "yes"

$ws.Range("B26").Value2 = "a slice of watermelon"
$ws.Range("C26").Value2 = "a slice of watermelon"
$ws.Range("D26").Value2 = "yes"

$ws.Range("B27").Value2 = "a slice of watermelon"
$ws.Range("C27").Value2 = "a slice of watermelon"
$ws.Range("D27").Value2 = "yes"

$ws.Range("B28").Value2 = "a slice of watermelon"
$ws.Range("C28").Value2 = "a slice of watermelon"
$ws.Range("D28").Value2 = "yes"

$ws.Range("B29").Value2 = "a slice of watermelon"
$ws.Range("C29").Value2 = "a slice of watermelon"
$ws.Range("D29").Value2 = "yes"

$ws.Range("B30").Value2 = "a slice of watermelon"
$ws.Range("C30").Value2 = "a slice of watermelon"
$ws.Range("D30").Value2 = "yes"

$ws.Range("D31").Value2 = "no"
$ws.Range("D32").Value2 = "no"

$ws.Range("B33").Value2 = "a flower"
$ws.Range("C33").Value2 = "a flower"
$ws.Range("D33").Value2 = "no"

$ws.Range("D34").Value2 = "no"
$ws.Range("D35").Value2 = "no"
$ws.Range("D36").Value2 = "no"

$ws.Range("B37").Value2 = "a flower"
$ws.Range("C37").Value2 = "a flower"
$ws.Range("D37").Value2 = "no"

$ws.Range("H26").Select()
